# Updated cryptos list on Fri Mar  8 23:55:50 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row (rows 2-51) on the active sheet. Values are written as literal
# text (matching the original inline-string cells): Price strings that would
# otherwise be auto-parsed by Excel as a plain number are written with a
# leading apostrophe to force text, then the cell style is reset back to
# "Normal" so no stray number-format/quote-prefix style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.176.53'
$ws.Range("E2").Value = '  +1.93%  '

$ws.Range("D3").Value = '3.885.75'
$ws.Range("E3").Value = '  +0.65%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '''484.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.43%  '

$ws.Range("D6").Value = '''144.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.63%  '

$ws.Range("E7").Value = '  -1.23%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("E9").Value = '  -2.84%  '

$ws.Range("E10").Value = '  +4.27%  '

$ws.Range("E11").Value = '  +12.01%  '

$ws.Range("D12").Value = '''42.67'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.53%  '

$ws.Range("E13").Value = '  +2.23%  '

$ws.Range("D14").Value = '4.508.27'
$ws.Range("E14").Value = '  +0.58%  '

$ws.Range("D15").Value = '''14.61'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.96%  '

$ws.Range("D16").Value = '3.895.75'
$ws.Range("E16").Value = '  +0.92%  '

$ws.Range("E17").Value = '  -0.32%  '

$ws.Range("D18").Value = '''19.70'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.08%  '

$ws.Range("E19").Value = '  -2.89%  '

$ws.Range("D20").Value = '68.187.24'
$ws.Range("E20").Value = '  +1.47%  '

$ws.Range("D21").Value = '''436.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.90%  '

$ws.Range("D22").Value = '''3.41'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.81%  '

$ws.Range("D23").Value = '''14.65'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.17%  '

$ws.Range("D24").Value = '''88.30'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.09%  '

$ws.Range("D25").Value = '''11.45'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +15.10%  '

$ws.Range("D26").Value = '''3.58'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.47%  '

$ws.Range("E27").Value = '  +4.55%  '

$ws.Range("D28").Value = '''38.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.31%  '

$ws.Range("E29").Value = '  +4.78%  '

$ws.Range("D30").Value = '''697.27'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.31%  '

$ws.Range("D31").Value = '''13.38'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.76%  '

$ws.Range("E32").Value = '  -2.58%  '

$ws.Range("E33").Value = '  +2.84%  '

$ws.Range("D34").Value = '0.0₃0910'
$ws.Range("E34").Value = '  +33.75%  '

$ws.Range("D35").Value = '''41.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.18%  '

$ws.Range("D36").Value = '''59.74'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.11%  '

$ws.Range("D37").Value = '''5.73'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.77%  '

$ws.Range("E38").Value = '  -6.69%  '

$ws.Range("D39").Value = '''0.999'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.00%  '

$ws.Range("E40").Value = '  -2.12%  '

$ws.Range("D41").Value = '''3.04'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.90%  '

$ws.Range("D42").Value = '''2.75'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.39%  '

$ws.Range("D43").Value = '''2.99'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.59%  '

$ws.Range("E44").Value = '  -1.48%  '

$ws.Range("D45").Value = '''0.142'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.31%  '

$ws.Range("E46").Value = '  -0.21%  '

$ws.Range("E47").Value = '  -1.46%  '

$ws.Range("E48").Value = '  -1.49%  '

$ws.Range("D49").Value = '''146.15'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.97%  '

$ws.Range("E50").Value = '  -2.29%  '

$ws.Range("D51").Value = '''2.84'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.64%  '
